$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2328912624424504
$ws.Range("D2").Value = 0.3174028231195223
$ws.Range("E2").Value = 0.0907717726627979
$ws.Range("F2").Value = 8.151866876652718
$ws.Range("G2").Value = 0.002693740383450863
$ws.Range("I2").Value = 0.5626248368130788
$ws.Range("J2").Value = 0.07877150001412048
$ws.Range("M2").Value = 4.497324019383484

$ws.Range("B3").Value = 0.2039938402802761
$ws.Range("D3").Value = 0.2841876498375768
$ws.Range("E3").Value = 0.0789721345030614
$ws.Range("F3").Value = 8.006350376761048
$ws.Range("G3").Value = 0.002711687236060236
$ws.Range("I3").Value = 0.5527629649614383
$ws.Range("J3").Value = 0.07639903279849136
$ws.Range("M3").Value = 4.17769716926

$ws.Range("B4").Value = 0.1862001470456391
$ws.Range("D4").Value = 0.2641461513723868
$ws.Range("E4").Value = 0.07173831418496235
$ws.Range("F4").Value = 7.926547348246714
$ws.Range("G4").Value = 0.002723206857165268
$ws.Range("I4").Value = 0.5472579255766874
$ws.Range("J4").Value = 0.07496044241923983
$ws.Range("M4").Value = 3.98359645714487

$ws.Range("B5").Value = 0.1789365688381395
$ws.Range("D5").Value = 0.2560619924672949
$ws.Range("E5").Value = 0.06879229350154503
$ws.Range("F5").Value = 7.896365601515896
$ws.Range("G5").Value = 0.002728027985388804
$ws.Range("I5").Value = 0.5451528675961228
$ws.Range("J5").Value = 0.07437878826331712
$ws.Range("M5").Value = 3.905017250361794

$ws.Range("B6").Value = 0.177729708878644
$ws.Range("D6").Value = 0.2547244665367145
$ws.Range("E6").Value = 0.06830319053859313
$ws.Range("F6").Value = 7.891493437547638
$ws.Range("G6").Value = 0.002728836217230569
$ws.Range("I6").Value = 0.5448116763684538
$ws.Range("J6").Value = 0.07428248274018046
$ws.Range("M6").Value = 3.891999904540114

$ws.Range("B7").Value = 0.1861022382144597
$ws.Range("D7").Value = 0.2640367973622517
$ws.Range("E7").Value = 0.07169857704098348
$ws.Range("F7").Value = 7.926130917207104
$ws.Range("G7").Value = 0.002723271361814451
$ws.Range("I7").Value = 0.5472289761188875
$ws.Range("J7").Value = 0.07495257943140388
$ws.Range("M7").Value = 3.982534640509272

$ws.Range("B8").Value = 0.2229380312659828
$ws.Range("D8").Value = 0.3058734614974981
$ws.Range("E8").Value = 0.08670031534071398
$ws.Range("F8").Value = 8.099674064923903
$ws.Range("G8").Value = 0.002699825264281586
$ws.Range("I8").Value = 0.5591102077466772
$ws.Range("J8").Value = 0.07794972871329264
$ws.Range("M8").Value = 4.386656101675925

$ws.Range("B9").Value = 0.2947678742724236
$ws.Range("D9").Value = 0.3909885525628169
$ws.Range("E9").Value = 0.1162558863173757
$ws.Range("F9").Value = 8.518627748919187
$ws.Range("G9").Value = 0.002657768992420295
$ws.Range("I9").Value = 0.5867839264241965
$ws.Range("J9").Value = 0.0839701482974391
$ws.Range("M9").Value = 5.197323384710359

$ws.Range("B10").Value = 0.3472971166443131
$ws.Range("D10").Value = 0.4557980457789199
$ws.Range("E10").Value = 0.1381260196927343
$ws.Range("F10").Value = 8.878602576530909
$ws.Range("G10").Value = 0.002629193096444872
$ws.Range("I10").Value = 0.609804863871986
$ws.Range("J10").Value = 0.088480382627111
$ws.Range("M10").Value = 5.805700348586214

$ws.Range("B11").Value = 0.3711427610551823
$ws.Range("D11").Value = 0.4858687095702408
$ws.Range("E11").Value = 0.1481259497588283
$ws.Range("F11").Value = 9.054672088856648
$ws.Range("G11").Value = 0.002616682340005267
$ws.Range("I11").Value = 0.6208682896732824
$ws.Range("J11").Value = 0.09055119052499805
$ws.Range("M11").Value = 6.085642332858953

$ws.Range("B12").Value = 0.3801653520100103
$ws.Range("D12").Value = 0.497348326357212
$ws.Range("E12").Value = 0.1519214464258383
$ws.Range("F12").Value = 9.123200786344285
$ws.Range("G12").Value = 0.002612013833554634
$ws.Range("I12").Value = 0.6251433085949856
$ws.Range("J12").Value = 0.09133809723518027
$ws.Range("M12").Value = 6.192141625164055

$ws.Range("B13").Value = 0.3782224973815573
$ws.Range("D13").Value = 0.4948717412123074
$ws.Range("E13").Value = 0.1511036065636446
$ws.Range("F13").Value = 9.108358018614524
$ws.Range("G13").Value = 0.002613016228865991
$ws.Range("I13").Value = 0.6242187928083496
$ws.Range("J13").Value = 0.09116850113999675
$ws.Range("M13").Value = 6.169182699201997

$ws.Range("B14").Value = 0.3718852000603476
$ws.Range("D14").Value = 0.4868112439215224
$ws.Range("E14").Value = 0.148438024129085
$ws.Range("F14").Value = 9.060272312142729
$ws.Range("G14").Value = 0.002616296882386284
$ws.Range("I14").Value = 0.6212182802911173
$ws.Range("J14").Value = 0.09061587498924695
$ws.Range("M14").Value = 6.094394050968333

$ws.Range("B15").Value = 0.3680024800529509
$ws.Range("D15").Value = 0.4818862449245103
$ws.Range("E15").Value = 0.1468064588274274
$ws.Range("F15").Value = 9.031062572760959
$ws.Range("G15").Value = 0.002618315332495386
$ws.Range("I15").Value = 0.6193915371226382
$ws.Range("J15").Value = 0.09027773178949872
$ws.Range("M15").Value = 6.048648841825212

$ws.Range("B16").Value = 0.3457377352266064
$ws.Range("D16").Value = 0.4538453640126932
$ws.Range("E16").Value = 0.1374736461960211
$ws.Range("F16").Value = 8.86735092340831
$ws.Range("G16").Value = 0.002630020432287695
$ws.Range("I16").Value = 0.6090937814918931
$ws.Range("J16").Value = 0.08834543369251691
$ws.Range("M16").Value = 5.787472401005289

$ws.Range("B17").Value = 0.3320661558088602
$ws.Range("D17").Value = 0.4367994036132359
$ws.Range("E17").Value = 0.131762321314099
$ws.Range("F17").Value = 8.770131832570996
$ws.Range("G17").Value = 0.002637325416988563
$ws.Range("I17").Value = 0.6029281913893385
$ws.Range("J17").Value = 0.087164911741759
$ws.Range("M17").Value = 5.628088453685422

$ws.Range("B18").Value = 0.3241979012330205
$ws.Range("D18").Value = 0.4270499150999854
$ws.Range("E18").Value = 0.1284820398347577
$ws.Range("F18").Value = 8.715367862508117
$ws.Range("G18").Value = 0.002641573111412901
$ws.Range("I18").Value = 0.5994375049925651
$ws.Range("J18").Value = 0.08648770273313744
$ws.Range("M18").Value = 5.536713029133182

$ws.Range("B19").Value = 0.3215330355147898
$ws.Range("D19").Value = 0.4237581007703568
$ws.Range("E19").Value = 0.1273721626568118
$ws.Range("F19").Value = 8.697021410071216
$ws.Range("G19").Value = 0.002643019254767343
$ws.Range("I19").Value = 0.5982651524636253
$ws.Range("J19").Value = 0.08625872015109337
$ws.Range("M19").Value = 5.505825094700697

$ws.Range("B20").Value = 0.3335220079601129
$ws.Range("D20").Value = 0.4386082377396576
$ws.Range("E20").Value = 0.1323698030988751
$ws.Range("F20").Value = 8.780360953064019
$ws.Range("G20").Value = 0.002636543030257661
$ws.Range("I20").Value = 0.60357877115554
$ws.Range("J20").Value = 0.08729039450540199
$ws.Range("M20").Value = 5.645024072790335

$ws.Range("B21").Value = 0.3737468142786895
$ws.Range("D21").Value = 0.4891762318374049
$ws.Range("E21").Value = 0.149220721606973
$ws.Range("F21").Value = 9.074345214122843
$ws.Range("G21").Value = 0.002615331410845236
$ws.Range("I21").Value = 0.6220972772708251
$ws.Range("J21").Value = 0.09077812044365174
$ws.Range("M21").Value = 6.116347682887664

$ws.Range("B22").Value = 0.3999938876883959
$ws.Range("D22").Value = 0.522768509064349
$ws.Range("E22").Value = 0.1602854193778995
$ws.Range("F22").Value = 9.277335420131408
$ws.Range("G22").Value = 0.002601870250399698
$ws.Range("I22").Value = 0.6346990967574442
$ws.Range("J22").Value = 0.09307350274676196
$ws.Range("M22").Value = 6.427265554285157

$ws.Range("B23").Value = 0.3859891847647248
$ws.Range("D23").Value = 0.5047872658294636
$ws.Range("E23").Value = 0.1543747751033067
$ws.Range("F23").Value = 9.16797361684786
$ws.Range("G23").Value = 0.002609018365313459
$ws.Range("I23").Value = 0.6279274165856918
$ws.Range("J23").Value = 0.09184695547372712
$ws.Range("M23").Value = 6.261047852967664

$ws.Range("B24").Value = 0.3328638429951241
$ws.Range("D24").Value = 0.4377903072483775
$ws.Range("E24").Value = 0.1320951506112138
$ws.Range("F24").Value = 8.775732858723359
$ws.Range("G24").Value = 0.002636896597993902
$ws.Range("I24").Value = 0.6032844759054399
$ws.Range("J24").Value = 0.08723365910987013
$ws.Range("M24").Value = 5.63736668282047

$ws.Range("B25").Value = 0.2753792222175662
$ws.Range("D25").Value = 0.3675932546972547
$ws.Range("E25").Value = 0.1082384150726767
$ws.Range("F25").Value = 8.39644603875422
$ws.Range("G25").Value = 0.002668733399738481
$ws.Range("I25").Value = 0.5788271299223453
$ws.Range("J25").Value = 0.08232618591369345
$ws.Range("M25").Value = 4.975900418004812
